$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared string text updates ---
$ws.Range("A8").Replace("46", "47") | Out-Null
$ws.Range("C9").Replace("11/13/2023", "11/20/2023") | Out-Null
$ws.Range("C9").Replace("11/19/2023", "11/26/2023") | Out-Null

# --- Numeric/text cell updates ---
$ws.Range("F14").Copy($ws.Range("G14"))
$ws.Range("E14").Copy($ws.Range("H14"))
$ws.Range("D15").Copy($ws.Range("C15"))
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 100
$ws.Range("L15").Value = -18.75
$ws.Range("N15").Value = -78.333333333333
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = 14
$ws.Range("H16").Value = 21.428571428571
$ws.Range("I16").Value = 178
$ws.Range("J16").Value = 176
$ws.Range("K16").Value = 1.136363636363
$ws.Range("L16").Value = 5.325443786982
$ws.Range("M16").Value = -22.608695652173
$ws.Range("N16").Value = -78.424242424242
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = -20
$ws.Range("F17").Value = 36
$ws.Range("G17").Value = 35
$ws.Range("H17").Value = 2.857142857142
$ws.Range("I17").Value = 390
$ws.Range("J17").Value = 359
$ws.Range("K17").Value = 8.635097493036
$ws.Range("L17").Value = 15.727002967359
$ws.Range("M17").Value = 88.405797101449
$ws.Range("N17").Value = -27.643784786641
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = -37.5
$ws.Range("F18").Value = 13
$ws.Range("H18").Value = -35
$ws.Range("I18").Value = 123
$ws.Range("J18").Value = 183
$ws.Range("K18").Value = -32.786885245901
$ws.Range("L18").Value = 29.473684210526
$ws.Range("M18").Value = 41.379310344827
$ws.Range("N18").Value = -71.12676056338
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 33
$ws.Range("G19").Value = 34
$ws.Range("H19").Value = -2.941176470588
$ws.Range("I19").Value = 392
$ws.Range("J19").Value = 439
$ws.Range("K19").Value = -10.706150341685
$ws.Range("L19").Value = 17.717717717717
$ws.Range("M19").Value = 65.400843881856
$ws.Range("N19").Value = -25.333333333333
$ws.Range("C20").Value = 2
$ws.Range("E20").Value = -33.333333333333
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = 12.5
$ws.Range("I20").Value = 77
$ws.Range("J20").Value = 80
$ws.Range("K20").Value = -3.75
$ws.Range("L20").Value = 42.592592592592
$ws.Range("M20").Value = 120
$ws.Range("N20").Value = -79.947916666666
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 28
$ws.Range("E21").Value = -21.428571428571
$ws.Range("F21").Value = 110
$ws.Range("G21").Value = 112
$ws.Range("H21").Value = -1.785714285714
$ws.Range("I21").Value = 1179
$ws.Range("J21").Value = 1259
$ws.Range("K21").Value = -6.354249404289
$ws.Range("L21").Value = 16.501976284585
$ws.Range("M21").Value = 43.25637910085
$ws.Range("N21").Value = -57.666068222621
$ws.Range("C23").Value = 6
$ws.Range("D23").Value = 4
$ws.Range("E23").Value = 50
$ws.Range("F23").Value = 40
$ws.Range("G23").Value = 30
$ws.Range("H23").Value = 33.333333333333
$ws.Range("I23").Value = 383
$ws.Range("J23").Value = 375
$ws.Range("K23").Value = 2.133333333333
$ws.Range("L23").Value = -2.295918367346
$ws.Range("M23").Value = 45.075757575757
$ws.Range("C24").Value = 13
$ws.Range("D24").Value = 16
$ws.Range("E24").Value = -18.75
$ws.Range("F24").Value = 76
$ws.Range("G24").Value = 72
$ws.Range("H24").Value = 5.555555555555
$ws.Range("I24").Value = 844
$ws.Range("J24").Value = 817
$ws.Range("K24").Value = 3.304773561811
$ws.Range("L24").Value = 23.032069970845
$ws.Range("M24").Value = 38.360655737704
$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 133.333333333333
$ws.Range("F25").Value = 71
$ws.Range("G25").Value = 37
$ws.Range("H25").Value = 91.891891891891
$ws.Range("I25").Value = 606
$ws.Range("J25").Value = 497
$ws.Range("K25").Value = 21.931589537223
$ws.Range("L25").Value = 22.672064777327
$ws.Range("M25").Value = -1.782820097244
$ws.Range("D26").Copy($ws.Range("C26"))
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 0
$ws.Range("L26").Value = 8.695652173913
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 1
$ws.Range("I27").Value = 52
$ws.Range("J27").Value = 61
$ws.Range("K27").Value = -14.754098360655
$ws.Range("L27").Value = -28.767123287671
$ws.Range("C28").Copy($ws.Range("F28"))
$ws.Range("C28").Copy($ws.Range("G28"))
$ws.Range("E28").Copy($ws.Range("H28"))
$ws.Range("M28").Value = -60.526315789473
$ws.Range("N28").Value = -83.870967741935
$ws.Range("C29").Copy($ws.Range("F29"))
$ws.Range("C29").Copy($ws.Range("G29"))
$ws.Range("E29").Copy($ws.Range("H29"))
$ws.Range("M29").Value = -57.575757575757
$ws.Range("N29").Value = -83.908045977011
